$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; this shifts all the existing rows
# (27..90) down to (28..91), preserving their values/formatting.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 44953
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108002
$ws.Range("J27").Value = "Mango"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 60
$ws.Range("N27").Value = 6500
$ws.Range("O27").Value = 7000
$ws.Range("P27").Value = 6750
$ws.Range("Q27").Value = "$/bandeja 4 kilos"
$ws.Range("R27").Value = "Perú"
$ws.Range("S27").Value = 1688
$ws.Range("T27").Value = 4
